# Apply the "cn181107" console upload edits to the elastic-network-interface
# delete confirmation dialog workbook.
#
# Summary of changes:
#  - The English confirmation message in column C (row 2) is reworded from
#      Confirm to delete elastic network interface"
#    to
#      Confirm to delete this elastic network interface "
#    (the trailing character is a left curly quote, matching the Chinese
#    source string that already ends with the same quote mark).
#  - The active selection / scroll position of the sheet is updated.
#  - Column C is widened so the longer English string is readable.
#  - Basic print/page setup (A4, portrait) is applied to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the confirmation message text (column C, row 2) ----------------
$ws.Range("C2").Value = "Confirm to delete this elastic network interface “"

# --- View / selection state -------------------------------------------
# Scroll the sheet so column D is the left-most visible column and select F14
# (matches the saved sheetView in the authored workbook).
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F14").Select()

# --- Column width -----------------------------------------------------
# Widen column C to fit the longer translated text.
$ws.Columns.Item(3).ColumnWidth = 34.15

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait
